$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Text number format on Price cells whose new values would otherwise
# be auto-parsed as numbers/dates by Excel, so they stay text like the source data.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"

# Apply the updated cell values
$ws.Range("D2").Value = "64.812.19"
$ws.Range("E2").Value = "  -0.71%  "
$ws.Range("D3").Value = "3.441.84"
$ws.Range("E3").Value = "  -1.22%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "573.74"
$ws.Range("E5").Value = "  -1.04%  "
$ws.Range("D6").Value = "159.12"
$ws.Range("E6").Value = "  -0.79%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("D8").Value = "3.441.93"
$ws.Range("E8").Value = "  -1.26%  "
$ws.Range("D9").Value = "0.582"
$ws.Range("E9").Value = "  -3.87%  "
$ws.Range("E10").Value = "  -1.43%  "
$ws.Range("E11").Value = "  -2.66%  "
$ws.Range("E12").Value = "  -1.21%  "
$ws.Range("D13").Value = "4.033.42"
$ws.Range("E13").Value = "  -1.34%  "
$ws.Range("E14").Value = "  -0.82%  "
$ws.Range("D15").Value = "27.66"
$ws.Range("E15").Value = "  -3.49%  "
$ws.Range("E16").Value = "  -5.74%  "
$ws.Range("D17").Value = "64.809.54"
$ws.Range("E17").Value = "  -0.79%  "
$ws.Range("D18").Value = "3.442.89"
$ws.Range("E18").Value = "  -0.80%  "
$ws.Range("E19").Value = "  -2.20%  "
$ws.Range("D20").Value = "13.90"
$ws.Range("E20").Value = "  -3.19%  "
$ws.Range("D21").Value = "381.55"
$ws.Range("E21").Value = "  -1.44%  "
$ws.Range("D22").Value = "7.97"
$ws.Range("E22").Value = "  -3.65%  "
$ws.Range("D23").Value = "0.548"
$ws.Range("E23").Value = "  -1.06%  "
$ws.Range("D24").Value = "0.997"
$ws.Range("E24").Value = "  -0.08%  "
$ws.Range("D25").Value = "72.15"
$ws.Range("E25").Value = "  -1.54%  "
$ws.Range("D26").Value = "0.0000119"
$ws.Range("E26").Value = "  -3.54%  "
$ws.Range("D27").Value = "9.83"
$ws.Range("E27").Value = "  -2.02%  "
$ws.Range("E28").Value = "  -0.97%  "
$ws.Range("E29").Value = "  -0.07%  "
$ws.Range("E30").Value = "  +0.05%  "
$ws.Range("D31").Value = "6.10"
$ws.Range("E31").Value = "  -2.83%  "
$ws.Range("E32").Value = "  -2.42%  "
$ws.Range("D33").Value = "23.25"
$ws.Range("E33").Value = "  -1.68%  "
$ws.Range("E34").Value = "  -3.42%  "
$ws.Range("D35").Value = "1.57"
$ws.Range("E35").Value = "  -0.69%  "
$ws.Range("D36").Value = "161.27"
$ws.Range("E36").Value = "  -0.81%  "
$ws.Range("E37").Value = "  -2.00%  "
$ws.Range("D38").Value = "2.900.18"
$ws.Range("E38").Value = "  -3.32%  "
$ws.Range("E39").Value = "  -4.21%  "
$ws.Range("D40").Value = "6.68"
$ws.Range("E40").Value = "  +2.14%  "
$ws.Range("D41").Value = "26.22"
$ws.Range("E41").Value = "  -3.48%  "
$ws.Range("E42").Value = "  -1.03%  "
$ws.Range("D43").Value = "42.92"
$ws.Range("E43").Value = "  +0.62%  "
$ws.Range("D44").Value = "0.0317"
$ws.Range("E44").Value = "  -2.49%  "
$ws.Range("D45").Value = "0.778"
$ws.Range("E45").Value = "  -0.07%  "
$ws.Range("D46").Value = "25.97"
$ws.Range("E46").Value = "  +1.62%  "
$ws.Range("B47").Value = "dogwifhat"
$ws.Range("C47").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D47").Value = "2.27"
$ws.Range("E47").Value = "  +2.78%  "
$ws.Range("B48").Value = "ONDO"
$ws.Range("C48").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D48").Value = "1.09"
$ws.Range("E48").Value = "  -2.46%  "
$ws.Range("D49").Value = "316.19"
$ws.Range("E49").Value = "  -2.02%  "
$ws.Range("E50").Value = "  -3.63%  "
$ws.Range("E51").Value = "  -3.51%  "
